$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D and E columns to Text format temporarily to prevent Excel
# from auto-converting numeric-looking strings (e.g. "1.00", "610.19")
# into actual numbers and losing formatting/precision.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.986.52"
$ws.Range("E2").Value = "  -3.79%  "
$ws.Range("D3").Value = "3.517.84"
$ws.Range("E3").Value = "  -4.37%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "610.19"
$ws.Range("E5").Value = "  -5.61%  "
$ws.Range("D6").Value = "151.97"
$ws.Range("E6").Value = "  -4.89%  "
$ws.Range("D7").Value = "3.516.15"
$ws.Range("E7").Value = "  -4.36%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("E10").Value = "  -4.35%  "
$ws.Range("D11").Value = "6.89"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("E13").Value = "  -5.26%  "
$ws.Range("D14").Value = "4.113.96"
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").Value = "31.60"
$ws.Range("E15").Value = "  -3.32%  "
$ws.Range("D16").Value = "3.514.35"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "66.988.88"
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -3.39%  "
$ws.Range("E20").Value = "  -4.35%  "
$ws.Range("D21").Value = "444.01"
$ws.Range("E21").Value = "  -5.32%  "
$ws.Range("D22").Value = "9.12"
$ws.Range("E22").Value = "  -9.62%  "
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").Value = "77.65"
$ws.Range("E24").Value = "  -2.51%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "3.659.20"
$ws.Range("E26").Value = "  -4.31%  "
$ws.Range("E27").Value = "  -2.97%  "
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  -7.15%  "
$ws.Range("D29").Value = "8.14"
$ws.Range("E29").Value = "  -10.44%  "
$ws.Range("E30").Value = "  -4.30%  "
$ws.Range("D31").Value = "1.64"
$ws.Range("E31").Value = "  -4.29%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "0.160"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "25.67"
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").Value = "6.15"
$ws.Range("E35").Value = "  -4.65%  "
$ws.Range("E36").Value = "  -7.23%  "
$ws.Range("D37").Value = "3.508.21"
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("D38").Value = "8.06"
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "173.14"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -5.55%  "
$ws.Range("D44").Value = "0.0860"
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("D45").Value = "0.888"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("D46").Value = "45.20"
$ws.Range("E46").Value = "  -4.28%  "
$ws.Range("D47").Value = "27.06"
$ws.Range("E47").Value = "  -6.26%  "
$ws.Range("E48").Value = "  -5.92%  "
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  -3.52%  "
$ws.Range("E51").Value = "  -5.19%  "

# Restore default cell style (removes the temporary text format,
# matching the original workbook which had no explicit style on these cells).
$ws.Range("D2:E51").Style = "Normal"

